$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8831.355113164813
$ws.Range("B5").Value = 79354.61403309148
$ws.Range("B7").Value = 1421.201574319291
$ws.Range("B9").Value = 2992.531864811288
$ws.Range("B10").Value = 159967.7468449901
$ws.Range("B11").Value = 0.08665150310063811
$ws.Range("B12").Value = 0.2631454521573777
$ws.Range("B13").Value = 0.3499999999999938
$ws.Range("B14").Value = 0.9980618842702891
$ws.Range("B15").Value = 0.8838669379974597
